$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 114, shifting existing rows 114-117 down to 115-118
$ws.Rows.Item(114).Insert()

# Populate the new row 114 with the new weekly price entry
$ws.Range("A114").Value = 10
$ws.Range("B114").Value = "Vega Modelo de Temuco"
$ws.Range("C114").Value = "La Araucanía"
$ws.Range("D114").Value = 44509
$ws.Range("E114").Value = 9
$ws.Range("F114").Value = 100114007
$ws.Range("G114").Value = "Jengibre"
$ws.Range("H114").Value = "Sin especificar"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 30
$ws.Range("K114").Value = 20000
$ws.Range("L114").Value = 20000
$ws.Range("M114").Value = 20000
$ws.Range("N114").Value = "$/caja 13 kilos"
$ws.Range("O114").Value = "Perú"
$ws.Range("P114").Value = 1538
$ws.Range("Q114").Value = 13
$ws.Range("R114").Value = "Hortaliza"
